# Populate the empty data cells of the REPRODUCIBILITY table (Lot 1-4,
# Standard Deviation, Mean) for Sample 1, Sample 2 and Sample 3 rows.
#
# The table is the last (7th) table in the document:
#   Header row: (blank) | Lot 1 | Lot 2 | Lot 3 | Lot 4 | Standard Deviation | Mean
#   Row 2: Sample 1 | 150  | 154  | 170  | 150  | 156  | 5.2%
#   Row 3: Sample 2 | 602  | 649  | 645  | 637  | 633  | 2.9%
#   Row 4: Sample 3 | 1476 | 1672 | 1722 | 1744 | 1654 | 7.2%

$d = $word.ActiveDocument

# Locate the REPRODUCIBILITY table by checking which table contains a
# "Lot 1" header cell, rather than hard-coding an index.
$table = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Cell(1, 2).Range.Text -like "Lot 1*") {
        $table = $candidate
        break
    }
}

$values = @(
    @("150", "154", "170", "150", "156", "5.2%"),
    @("602", "649", "645", "637", "633", "2.9%"),
    @("1476", "1672", "1722", "1744", "1654", "7.2%")
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $rowIndex = $r + 2
    $rowValues = $values[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $colIndex = $c + 2
        $table.Cell($rowIndex, $colIndex).Range.Text = $rowValues[$c]
    }
}
